# edit.ps1 - applies the transcription corrections described by the commit
# "ajuste bo histórico na parte do texto e inserido documentação"
#
# Strategy: use Find/Execute (wildcards off, MatchCase on) to do precise,
# surgical replacements of the exact old phrasing with the new phrasing.
# A handful of edits also re-flow paragraph breaks (splitting one sentence
# into two, or merging two paragraphs into one / relabeling a speaker), so
# those use "^p" inside the search/replacement text to span paragraph marks.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $found = $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Output "NOT FOUND: $old"
    }
}

# 1) "Cliente" — conta atrasada / 104 reais paragraph: two sentence splits.
Replace-Text "Aí eu paguei 104 reais esses dias e a minha internet está lenta" "Aí eu paguei 104 reais esses dias. E a minha internet está lenta"
Replace-Text "Eu não sei se liberaram, não sei se está ativa, eu não sei." "Eu não sei se liberaram. Não sei se está ativa, eu não sei."

# 2) "Atendente" — equipamento além do nosso.
Replace-Text "Eu vou verificar pra senhora. Tem algum tipo de equipamento além do nosso, da Leste?" "Eu vou verificar pra senhora. A senhora tem algum tipo de equipamento além do nosso da Leste?"

# 3) "Cliente" — dois aparelhos.
Replace-Text "Eu tenho dois aparelhos. Tem o roteador, tem o roteador e tem o outro, não, aquele tem uma anteninha, um cabinho assim. A, tem um roteador e tem o outro." "Eu tenho dois aparelhos. E tem o roteador. Tem o roteador e tem o outro, não, aquele tem uma anteninha, um cabinho assim. A, tem um roteador e tem o outrozinho."

# 4) "Cliente" — Não, quando eu entro aqui / Não consigo entrar porque trava.
Replace-Text "**Cliente:** Não, quando eu entro aqui, tipo assim," "**Cliente:** Não. Quando eu entro aqui, tipo assim,"
Replace-Text "Não consigo entrar porque trava," "Não consigo entrar, porque trava,"

# 5) "Atendente" — Compreendo, senhora Kelly / gerência nas configurações.
Replace-Text "a gente não tem ger... nas configura... por não ser um equipamento" "a gente não tem gerência nas configurações, por não ser um equipamento"

# 6) "Cliente" — Eu consigo... pera aí... branco.
Replace-Text "o meu é tranquilo porque eu tenho um... preço... um branco, tenho que ter um cabelo amarelo." "o meu é tranquilo porque eu tenho um, é branco, tenho que ter um cabelo amarelo."

# 7) "Cliente" — Não, porque já me ofereceram esse plano.
Replace-Text "Quando ventava, eu não ficava sem internet porque lá onde eu morava ventava muito." "Quando ventava, eu não ficava sem internet. Porque lá onde eu morava ventava muito."
Replace-Text "Eu não sei se isso que acontece porque diminui o mega na hora da fatura, não sei. Se não foi a paga, não diz certo." "Eu não sei se isso que acontece porque diminui o mega na hora da fatura, não? Se não foi a paga, não diz certo."

# 8) "Cliente" — Agora eu estou começando a ter.
Replace-Text "A internet por aí... Estão de onde para poder... A gente trocar de plano... Como é que foi minha primeira hora?" "A internet por aí. Estão de onde para poder... A gente trocar de plano. Como é que foi minha primeira hora?"

# 9) Merge the "Atendente" (Leste Suporte) paragraph with the following
#    "Cliente" one-liner ("Mas a minha fatura continua sendo de 100 reais?")
#    into a single Atendente paragraph (also split a comma off as a new
#    sentence), removing the now-empty paragraph gap between them.
Replace-Text "Se a senhora tiver disponibilidade de baixar o aplicativo Leste Suporte, aí eu vou baixar isso aí.^p^p**Cliente:** Mas a minha fatura continua sendo de 100 reais?" "Se a senhora tiver disponibilidade de baixar o aplicativo Leste Suporte. Aí eu vou baixar isso aí. Mas a minha fatura continua sendo de 100 reais?"

# 10) Re-label the next paragraph (formerly "Atendente", now logically the
#     client's reply) from Atendente to Cliente.
Replace-Text "**Atendente:** Não, não vou mandar o plano da senhora. Eu vou verificar as configurações do roteador da senhora." "**Cliente:** Não, não vou mandar o plano da senhora. Eu vou verificar as configurações do roteador da senhora."

# 11) Merge the short "Cliente: Ah, sim." paragraph into the following
#     "Atendente" paragraph, dropping the now-empty paragraph between them.
Replace-Text "**Cliente:** Ah, sim.^p^p**Atendente:** Aí na lojinha de aplicativos do celular" "**Atendente:** Ah, sim. Aí na lojinha de aplicativos do celular"
